# Refresh the crypto price/volume table (GitHub Actions scheduled update).
#
# Most "Price" (column D) / "Volume(1h)" (column E) cells are plain text
# (e.g. "42.760.41", "  +0.16%  ") and can just be re-assigned via .Value.
#
# A handful of new Price values parse as plain decimals (e.g. "256.72",
# "8.89"); assigning those directly would make Excel auto-convert the cell
# to a real number. To keep them as literal text (matching the original
# file, and avoiding a float-precision/style diff), we briefly mark the
# cell as Text format, assign the string, then reset the cell style back
# to Normal/default so no stray number-format style lingers behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.753.37"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.251.84"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "295.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("E14").Value = "  +23.47%  "
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "2.589.48"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "2.248.68"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "42.718.96"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +17.67%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "256.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.22%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.05%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.62%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.66%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.60%  "
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.89%  "
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("E41").Value = "  -5.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("E45").Value = "  -5.24%  "
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.15%  "
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.28%  "
